$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.049.04"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "'2.498.58"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'319.54"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").Value = "'105.74"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -4.29%  "
$ws.Range("D10").Value = "'38.70"
$ws.Range("E10").Value = "  -4.25%  "
$ws.Range("D11").Value = "'19.99"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "'2.890.92"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "'2.499.03"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").Value = "'0.832"
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").Value = "'47.901.57"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "'13.03"
$ws.Range("E19").Value = "  -2.91%  "
$ws.Range("D20").Value = "'2.96"
$ws.Range("E20").Value = "  +7.91%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'0.0₃0933"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").Value = "'71.18"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("D24").Value = "'272.08"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "'2.52"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'25.75"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E28").Value = "  +9.46%  "
$ws.Range("D29").Value = "'9.73"
$ws.Range("E29").Value = "  -4.73%  "
$ws.Range("D30").Value = "'0.140"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "'34.78"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("D32").Value = "'49.11"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "'19.13"
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").Value = "'0.0773"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").Value = "'1.93"
$ws.Range("E37").Value = "  -2.95%  "
$ws.Range("D38").Value = "'4.56"
$ws.Range("E38").Value = "  -3.95%  "
$ws.Range("E39").Value = "  -4.33%  "
$ws.Range("D40").Value = "'122.64"
$ws.Range("E40").Value = "  +2.65%  "
$ws.Range("D41").Value = "'0.110"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").Value = "'22.17"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").Value = "'2.001.62"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'3.15"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").Value = "'1.89"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'2.01"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").Value = "'8.88"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").Value = "'78.76"
$ws.Range("E51").Value = "  -2.41%  "
